$wb = $excel.ActiveWorkbook

# Sheet "studies" -> add new "PMID" column (H) header
$studies = $wb.Worksheets.Item("studies")
$studies.Activate()
$studies.Range("H1").Value = "PMID"
$studies.Range("H1").Select()

# Sheet "counts" -> add new "notes" column (F) header
$counts = $wb.Worksheets.Item("counts")
$counts.Activate()
$counts.Range("F1").Value = "notes"

# "counts" ends up as the active sheet/tab, with F2 selected
$counts.Range("F2").Select()
